# Add Q3-2022 data:
#  1. Insert a new "2022-Q3" summary row at the top of the "总计" (Total) sheet.
#  2. Insert a new "2022-Q3" worksheet (fund holdings detail) right after "总计"
#     and before "2022-Q2", pushing every later quarter sheet one slot to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" summary sheet: insert new row 2 and fill it in, then renumber the
#    running index in column A for the rows pushed down.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Restore the header-row-2 formatting (border/alignment) on the freshly
# inserted row by cloning it from the row directly below (old row 2, now row 3).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 22
$summary.Range("D2").Value = 1.31

# The rows that shifted down keep their original 0-based index; bump them by one.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet with the fund-holding detail table, placed right
#    after "总计".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$rows = @(
    @(0,  "000478", "建信中证500指数增强A",           "45.95", "82.53", "0.92", "0.4227", 5),
    @(1,  "006972", "金鹰民安回报一年定期开放混合A",   "26.02", "22.96", "1.62", "0.4215", 3),
    @(2,  "015110", "惠升领先优选混合A",               "1.50",  "81.38", "4.86", "0.0729", 4),
    @(3,  "014786", "惠升品质优选混合A",               "1.51",  "82.37", "4.81", "0.0726", 4),
    @(4,  "007735", "金鹰民安回报一年定期开放混合C",   "2.59",  "22.96", "1.62", "0.0420", 3),
    @(5,  "011351", "金鹰年年邮益一年持有期混合A",     "3.43",  "34.33", "1.11", "0.0381", 4),
    @(6,  "015043", "西部利得时代动力混合A",           "0.68",  "79.48", "5.53", "0.0376", 2),
    @(7,  "562500", "华夏中证机器人ETF",               "1.57",  "99.51", "2.07", "0.0325", 10),
    @(8,  "005633", "建信中证500指数增强C",           "3.42",  "82.53", "0.92", "0.0315", 5),
    @(9,  "006502", "财通集成电路产业股票A",           "0.72",  "83.09", "4.29", "0.0309", 8),
    @(10, "159770", "天弘中证机器人ETF",               "1.04",  "99.68", "2.07", "0.0215", 10),
    @(11, "006503", "财通集成电路产业股票C",           "0.45",  "83.09", "4.29", "0.0193", 8),
    @(12, "562360", "银华中证机器人ETF",               "0.77",  "97.23", "2.03", "0.0156", 10),
    @(13, "015097", "东财数字经济优选混合C",           "0.39",  "89.25", "3.78", "0.0147", 10),
    @(14, "015096", "东财数字经济优选混合A",           "0.31",  "89.25", "3.78", "0.0117", 10),
    @(15, "015044", "西部利得时代动力混合C",           "0.16",  "79.48", "5.53", "0.0088", 2),
    @(16, "002564", "新沃通盈灵活配置混合",           "0.11",  "93.51", "3.71", "0.0041", 10),
    @(17, "010466", "鹏扬景创混合C",                   "0.29",  "32.63", "1.06", "0.0031", 10),
    @(18, "011352", "金鹰年年邮益一年持有期混合C",     "0.27",  "34.33", "1.11", "0.0030", 4),
    @(19, "010465", "鹏扬景创混合A",                   "0.21",  "32.63", "1.06", "0.0022", 10),
    @(20, "015111", "惠升领先优选混合C",               "0.00",  "81.38", "4.86", 0,        4),
    @(21, "014787", "惠升品质优选混合C",               "0.00",  "82.37", "4.81", 0,        4)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}
